$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.187.43'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.262.21'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.580'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0840'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.64'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '2.610.30'
$ws.Range("E14").Value = '  +1.75%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.61'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.861'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '2.270.52'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '44.116.35'
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").Value = '0.0₃0988'
$ws.Range("E20").Value = '  +2.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.49%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.78%  '
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0851'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.71%  '
$ws.Range("E36").Value = '  +9.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.72%  '
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0317'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '1.781.78'
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '82.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '74.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("E51").Value = '  +4.54%  '
